# Auto-update draw results: append the 2025-11-19 Pick 4 row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

$row = 64
$rng = $ws.Range("A" + $row + ":E" + $row)

# Pre-format the new row as Text so date-/number-looking strings
# (the date, the 6-digit phase code) round-trip as literal text,
# matching every other row in the sheet.
$rng.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-11-19"
$ws.Cells.Item($row, 2).Value = "Pick 4"
$ws.Cells.Item($row, 3).Value = "251119"
$ws.Cells.Item($row, 4).Value = "9-7-2-5"
$ws.Cells.Item($row, 5).Value = "2025-11-19T21:37:21.821+04:00"

# Drop back to the workbook's default ("Normal") cell style so the new
# row doesn't carry a distinct style from the rest of the table.
$rng.Style = "Normal"
